$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 170.53847
$ws.Range("I55").Value = 173.7
$ws.Range("J55").Value = 160
$ws.Range("K55").Value = 173.7
$ws.Range("L55").Value = 160
$ws.Range("M55").Value = 40.30000000000001
$ws.Range("N55").Value = -588

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1700
$ws.Range("I70").Value = 1471.4286
$ws.Range("J70").Value = 1814.2858
$ws.Range("K70").Value = 4414.2858
$ws.Range("L70").Value = 5442.857400000001
$ws.Range("M70").Value = -4144.2858
$ws.Range("N70").Value = -5982.857400000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 1700
$ws.Range("I73").Value = 1471.4286
$ws.Range("J73").Value = 1814.2858
$ws.Range("K73").Value = 4414.2858
$ws.Range("L73").Value = 5442.857400000001
$ws.Range("M73").Value = -3478.2858
$ws.Range("N73").Value = -7314.857400000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 105979.27
$ws.Range("I138").Value = 1515.5294
$ws.Range("J138").Value = 127903.766
$ws.Range("K138").Value = 4546.5882
$ws.Range("L138").Value = 383711.298
$ws.Range("M138").Value = 593.4117999999999
$ws.Range("N138").Value = -393991.298

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24759.205
$ws.Range("I32").Value = 18522.275
$ws.Range("J32").Value = 43470
$ws.Range("K32").Value = 18522.275
$ws.Range("L32").Value = 43470
$ws.Range("M32").Value = -18235.275
$ws.Range("N32").Value = -44044

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2138.4666
$ws.Range("I61").Value = 1458.7368
$ws.Range("J61").Value = 3312.5454
$ws.Range("K61").Value = 1458.7368
$ws.Range("L61").Value = 3312.5454
$ws.Range("M61").Value = -1246.7368
$ws.Range("N61").Value = -3736.5454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 26771.41
$ws.Range("I74").Value = 28857.916
$ws.Range("J74").Value = 1733.3334
$ws.Range("K74").Value = 28857.916
$ws.Range("L74").Value = 1733.3334
$ws.Range("M74").Value = -27983.916
$ws.Range("N74").Value = -3481.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 26771.41
$ws.Range("I77").Value = 28857.916
$ws.Range("J77").Value = 1733.3334
$ws.Range("K77").Value = 144289.58
$ws.Range("L77").Value = 8666.666999999999
$ws.Range("M77").Value = -139921.58
$ws.Range("N77").Value = -17402.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2138.4666
$ws.Range("I136").Value = 1458.7368
$ws.Range("J136").Value = 3312.5454
$ws.Range("K136").Value = 4376.2104
$ws.Range("L136").Value = 9937.636200000001
$ws.Range("M136").Value = -1826.2104
$ws.Range("N136").Value = -15037.6362

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H27").Value = 28449.5
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 28449.5
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 28449.5
$ws.Range("N27").Value = -28833.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 33917.555
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 33917.555
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 33917.555
$ws.Range("N81").Value = -36039.555

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H84").Value = 33917.555
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 33917.555
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 101752.665
$ws.Range("N84").Value = -112360.665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 42305.6
$ws.Range("I31").Value = 41645.04
$ws.Range("J31").Value = 45608.4
$ws.Range("K31").Value = 41645.04
$ws.Range("L31").Value = 45608.4
$ws.Range("M31").Value = -41350.04
$ws.Range("N31").Value = -46198.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 42305.6
$ws.Range("I34").Value = 41645.04
$ws.Range("J34").Value = 45608.4
$ws.Range("K34").Value = 41645.04
$ws.Range("L34").Value = 45608.4
$ws.Range("M34").Value = -41443.04
$ws.Range("N34").Value = -46012.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3363.1428
$ws.Range("I99").Value = 3101.3333
$ws.Range("J99").Value = 3559.5
$ws.Range("K99").Value = 3101.3333
$ws.Range("L99").Value = 3559.5
$ws.Range("M99").Value = -1603.3333
$ws.Range("N99").Value = -6555.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3363.1428
$ws.Range("I126").Value = 3101.3333
$ws.Range("J126").Value = 3559.5
$ws.Range("K126").Value = 9303.999899999999
$ws.Range("L126").Value = 10678.5
$ws.Range("M126").Value = -6833.999899999999
$ws.Range("N126").Value = -15618.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 24083.355
$ws.Range("I132").Value = 27743.525
$ws.Range("J132").Value = 4213.857
$ws.Range("K132").Value = 83230.57500000001
$ws.Range("L132").Value = 12641.571
$ws.Range("M132").Value = -80700.57500000001
$ws.Range("N132").Value = -17701.571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 348537.38
$ws.Range("I4").Value = 918332.9399999999
$ws.Range("J4").Value = 329
$ws.Range("K4").Value = 2754998.82
$ws.Range("L4").Value = 987
$ws.Range("M4").Value = -2754886.82
$ws.Range("N4").Value = -1211

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1000002
$ws.Range("I68").Value = 1000002
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 3000006
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -2999195
$ws.Range("N68").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1000002
$ws.Range("I71").Value = 1000002
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 9000018
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -8995962
$ws.Range("N71").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 543.8570999999999
$ws.Range("I98").Value = 434.33334
$ws.Range("J98").Value = 626
$ws.Range("K98").Value = 1303.00002
$ws.Range("L98").Value = 1878
$ws.Range("M98").Value = 194.9999800000001
$ws.Range("N98").Value = -4874

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 186033.33
$ws.Range("I131").Value = 411.1111
$ws.Range("J131").Value = 223157.78
$ws.Range("K131").Value = 1233.3333
$ws.Range("L131").Value = 669473.34
$ws.Range("M131").Value = 3806.6667
$ws.Range("N131").Value = -679553.34

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 2150.8572
$ws.Range("I138").Value = 1139.2
$ws.Range("J138").Value = 3318.1538
$ws.Range("K138").Value = 3417.6
$ws.Range("L138").Value = 9954.4614
$ws.Range("M138").Value = 1722.4
$ws.Range("N138").Value = -20234.4614

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 3333.6572
$ws.Range("I140").Value = 1428.5333
$ws.Range("J140").Value = 4762.5
$ws.Range("K140").Value = 4285.5999
$ws.Range("L140").Value = 14287.5
$ws.Range("M140").Value = 894.4000999999998
$ws.Range("N140").Value = -24647.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 12000
$ws.Range("I75").Value = 12000
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 12000
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -11064
$ws.Range("N75").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H78").Value = 12000
$ws.Range("I78").Value = 12000
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 36000
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -31320
$ws.Range("N78").ClearContents()
